$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 (IP column) changes value from "192.168.1.113" to "127.0.0.1"
$ws.Range("C2").Value = "127.0.0.1"

# E2 (SqlIP column) keeps its text value "192.168.0.24" but now gets the
# same text-number-format as C2 (numFmtId 49 / "@")
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("E2").NumberFormat = $ws.Range("C2").NumberFormat

# Move the active cell selection from C2 to E2
$ws.Range("E2").Select()
